# Update the "nota_view" column (J): every cell currently valued 5 becomes 4.
# This reflects the commit updating attendance control (19/09/2022) and forum
# grades for the week 11/09/2022 - 17/09/2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range to know how many rows contain data.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # Column J = 10
    if ($cell.Value2 -eq 5) {
        $cell.Value = 4
    }
}
